# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 07:16"

# Row 5 - India
$ws.Range("B5").Value = 5118253
$ws.Range("C5").Value = 2360
$ws.Range("D5").Value = 4025079
$ws.Range("E5").Value = 1009944

# Row 27 - Israel
$ws.Range("B27").Value = 171768
$ws.Range("C27").Value = 1303
$ws.Range("D27").Value = 125619
$ws.Range("E27").Value = 44984

# Row 59 - Uzbekistan
$ws.Range("B59").Value = 49162
$ws.Range("C59").Value = 147
$ws.Range("D59").Value = 45474
$ws.Range("E59").Value = 3279
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 409

# Row 69 - Kenia
$ws.Range("B69").Value = 36393
$ws.Range("D69").Value = 23529
$ws.Range("E69").Value = 12227
$ws.Range("H69").Value = 637

# Row 76 - Australia
$ws.Range("B76").Value = 26813
$ws.Range("C76").Value = 34
$ws.Range("D76").Value = 23792
$ws.Range("E76").Value = 2189

# Row 131 - Tailandia
$ws.Range("D131").Value = 3325
$ws.Range("E131").Value = 107
